# Applies the "Added some exercises and wods" commit to Data.xlsx
$wb = $excel.ActiveWorkbook

$wsEquipment = $wb.Worksheets.Item("Equipment")
$wsExercises = $wb.Worksheets.Item("Exercises")
$wsWODs      = $wb.Worksheets.Item("WODs")
$wsBacklog   = $wb.Worksheets.Item("Backlog")

# ---------------------------------------------------------------------------
# Exercises sheet: new exercises (rows 6-21)
# Columns: A=Id  B=Name  C=ShortName(formula)  D=RequiredEquipment
#          E=TargetGroups  F=Difficulty  G=Description  H=Icon(formula)
# ---------------------------------------------------------------------------
$wsExercises.Range("A6").Value = 5
$wsExercises.Range("B6").Value = "Burpee"
$wsExercises.Range("E6").Value = "fullbody"
$wsExercises.Range("F6").Value = "Beginner"
$wsExercises.Range("G6").Value = "TODO"

$wsExercises.Range("A7").Value = 6
$wsExercises.Range("B7").Value = "Boxjump"
$wsExercises.Range("D7").Value = "[Box]"
$wsExercises.Range("E7").Value = "legs"
$wsExercises.Range("F7").Value = "Beginner"
$wsExercises.Range("G7").Value = "TODO"

$wsExercises.Range("A8").Value = 7
$wsExercises.Range("B8").Value = "Toes 2 bar"
$wsExercises.Range("D8").Value = "[pull-up bar]"
$wsExercises.Range("E8").Value = "core"
$wsExercises.Range("F8").Value = "Intermediate"
$wsExercises.Range("G8").Value = "TODO"

$wsExercises.Range("A9").Value = 8
$wsExercises.Range("B9").Value = "Wallbal"
$wsExercises.Range("D9").Value = "[medicine ball]"
$wsExercises.Range("E9").Value = "legs, shoulder, arms"
$wsExercises.Range("F9").Value = "Beginner"
$wsExercises.Range("G9").Value = "TODO"

$wsExercises.Range("A10").Value = 9
$wsExercises.Range("B10").Value = "Benchpress"
$wsExercises.Range("D10").Value = "[Bench], [weight lifting bar]"
$wsExercises.Range("F10").Value = "Intermediate"
$wsExercises.Range("G10").Value = "TODO"

$wsExercises.Range("A11").Value = 10
$wsExercises.Range("B11").Value = "Handstand pushup"
$wsExercises.Range("F11").Value = "Expert"
$wsExercises.Range("G11").Value = "TODO"

$wsExercises.Range("A12").Value = 11
$wsExercises.Range("B12").Value = "Goblet squat"
$wsExercises.Range("D12").Value = "[Kettlebell]"
$wsExercises.Range("E12").Value = "legs"
$wsExercises.Range("F12").Value = "Beginner"
$wsExercises.Range("G12").Value = "TODO"

$wsExercises.Range("A13").Value = 12
$wsExercises.Range("B13").Value = "Duble unders"
$wsExercises.Range("D13").Value = "[Jump rope]"
$wsExercises.Range("F13").Value = "Intermediate"
$wsExercises.Range("G13").Value = "TODO"

$wsExercises.Range("A14").Value = 12
$wsExercises.Range("B14").Value = "Singles"
$wsExercises.Range("D14").Value = "[Jump rope]"
$wsExercises.Range("F14").Value = "Beginner"
$wsExercises.Range("G14").Value = "TODO"

$wsExercises.Range("A15").Value = 13
$wsExercises.Range("B15").Value = "Russian swings"
$wsExercises.Range("D15").Value = "[Kettlebell]"
$wsExercises.Range("F15").Value = "Beginner"
$wsExercises.Range("G15").Value = "TODO"

$wsExercises.Range("A16").Value = 14
$wsExercises.Range("B16").Value = "Turkish getup"
$wsExercises.Range("E16").Value = "fullbody"
$wsExercises.Range("F16").Value = "Beginner"
$wsExercises.Range("G16").Value = "TODO"

$wsExercises.Range("A17").Value = 15
$wsExercises.Range("B17").Value = "Run"
$wsExercises.Range("E17").Value = "legs"
$wsExercises.Range("F17").Value = "Beginner"
$wsExercises.Range("G17").Value = "TODO"

$wsExercises.Range("A18").Value = 16
$wsExercises.Range("B18").Value = "Kettlebell swing"
$wsExercises.Range("D18").Value = "[Kettlebell]"
$wsExercises.Range("E18").Value = "core, arms, shoulder"
$wsExercises.Range("F18").Value = "Beginner"
$wsExercises.Range("G18").Value = "TODO"

$wsExercises.Range("A19").Value = 17
$wsExercises.Range("B19").Value = "Snatch"
$wsExercises.Range("D19").Value = "[weight lifting bar]"
$wsExercises.Range("F19").Value = "Expert"
$wsExercises.Range("G19").Value = "TODO"

$wsExercises.Range("A20").Value = 18
$wsExercises.Range("B20").Value = "Power snatch"
$wsExercises.Range("D20").Value = "[weight lifting bar]"
$wsExercises.Range("F20").Value = "Expert"
$wsExercises.Range("G20").Value = "TODO"

$wsExercises.Range("A21").Value = 19
$wsExercises.Range("B21").Value = "Hang power snatch"
$wsExercises.Range("D21").Value = "[weight lifting bar]"
$wsExercises.Range("F21").Value = "Intermediate"
$wsExercises.Range("G21").Value = "TODO"

# Exercises sheet also picked up an explicit (portrait / letter-ish A4-ish)
# page setup in the saved file.
$wsExercises.PageSetup.PaperSize = 9
$wsExercises.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# WODs sheet: new workouts (rows 5-10)
# Columns: A=Id  B=Name  C=ShortName  D=Type  E=Work  F=Rest  G=Rounds  H=Exercises(description)
# ---------------------------------------------------------------------------
$wsWODs.Range("A5").Value = 4
$wsWODs.Range("D5").Value = "Time"
$wsWODs.Range("G5").Value = 4
$wsWODs.Range("H5").Value = "For time: 10 [toes 2 bar], 15 [boxjumps], 20 [wallbals]"

$wsWODs.Range("A6").Value = 5
$wsWODs.Range("B6").Value = "Lynne"
$wsWODs.Range("D6").Value = "Rounds"
$wsWODs.Range("G6").Value = 5
$wsWODs.Range("H6").Value = "5 rounds not for time: Max reps [benchpress] @BW, Max reps strict [pullups]"

$wsWODs.Range("A7").Value = 6
$wsWODs.Range("D7").Value = "Amrap"
$wsWODs.Range("E7").Value = 120
$wsWODs.Range("F7").Value = 60
$wsWODs.Range("G7").Value = 3
$wsWODs.Range("H7").Value = "3 X Amrap 2: [Burpees], 60 sec rest between AMRAPs"

$wsWODs.Range("A8").Value = 7
$wsWODs.Range("D8").Value = "OTM"
$wsWODs.Range("G8").Value = 20
$wsWODs.Range("H8").Value = "On the minute: Odd: 10 [russian swings], 10 [double unders]. Even: 10 [goblet squats], 10 [double unders]"

$wsWODs.Range("A9").Value = 8
$wsWODs.Range("B9").Value = "Helen"
$wsWODs.Range("D9").Value = "Time"
$wsWODs.Range("H9").Value = "For time: 400 m [run], 21 [kettlebell swings], 12 [pullups]"

$wsWODs.Range("A10").Value = 9
$wsWODs.Range("D10").Value = "Time"
$wsWODs.Range("H10").Value = "For time. 21-15-9 [goblet squats], [burpees]"

# ---------------------------------------------------------------------------
# Backlog sheet: new backlog item (row 16)
# Columns: A=Id  B=Name  C=Description  D=Priority  E=Responsible
# ---------------------------------------------------------------------------
$wsBacklog.Range("A16").Value = 15
$wsBacklog.Range("B16").Value = "Registrering af resultater"
$wsBacklog.Range("D16").Value = "Medium"

# ---------------------------------------------------------------------------
# Selections / active sheet, matching the edited file's last-saved UI state.
# Order matters: the workbook was last saved with the WODs tab active.
# ---------------------------------------------------------------------------
$wsEquipment.Activate()
$wsEquipment.Range("A4").Select() | Out-Null

$wsExercises.Activate()
$wsExercises.Range("F22").Select() | Out-Null

$wsWODs.Activate()
$wsWODs.Range("A11").Select() | Out-Null

$wsBacklog.Activate()
$wsBacklog.Range("D16").Select() | Out-Null

$wsWODs.Activate()
